# Add a new worksheet "magapoke_2025-12-03" at the end of the workbook,
# replicating the layout/style of the existing weekly ranking sheets
# (header row "rank"/"title" in the existing bold+bordered+centered style,
# followed by 100 rows of rank/title data).

$wb = $excel.ActiveWorkbook

# Template sheet to copy the header formatting from (any existing weekly
# ranking sheet uses the same "rank"/"title" header style).
$templateSheet = $wb.Worksheets.Item(1)
$templateHeader = $templateSheet.Range("A1:B1")

# Insert the new sheet after the current last sheet, so it lands at the end
# of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "magapoke_2025-12-03"

# Header row.
$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"

# Copy the header's formatting (bold font, thin border, centered/top
# alignment) from the template sheet so the new sheet matches the others.
$templateHeader.Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

# Match the page margins used by the other weekly ranking sheets
# (left/right 0.75", top/bottom 1", header/footer 0.5").
$newSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$newSheet.PageSetup.TopMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$newSheet.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$newSheet.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Rank / title data rows (rows 2-101).
$titles = @(
    @(1, "ブルーロック"),
    @(2, "みいちゃんと山田さん"),
    @(3, "WIND BREAKER"),
    @(4, "信じていた仲間達にダンジョン奥地で殺されかけたがギフト『無限ガチャ』でレベル9999の仲間達を手に入れて元パーティーメンバーと世界に復讐＆『ざまぁ！』します！"),
    @(5, "東京卍リベンジャーズ"),
    @(6, "ガチアクタ"),
    @(7, "ギルティサークル"),
    @(8, "薫る花は凛と咲く"),
    @(9, "島耕作"),
    @(10, "イレギュラーズ"),
    @(11, "ベイビーステップ"),
    @(12, "十字架のろくにん"),
    @(13, "となりの黒川さん"),
    @(14, "黒猫と魔女の教室"),
    @(15, "転生したら第七王子だったので、気ままに魔術を極めます"),
    @(16, "南海トラフ巨大地震"),
    @(17, "愛妻の裏アカ"),
    @(18, "黄昏町プリズナーズ"),
    @(19, "魔女と傭兵"),
    @(20, "君が僕らを悪魔と呼んだ頃"),
    @(21, "味方が弱すぎて補助魔法に徹していた宮廷魔法師、追放されて最強を目指す"),
    @(22, "転生貴族、鑑定スキルで成り上がる～弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた～"),
    @(23, "蒼く染めろ"),
    @(24, "K-9~警視庁公安部公安第9課異能対策係~"),
    @(25, "異世界ウォーキング"),
    @(26, "魔術ギルド総帥～生まれ変わって今更やり直す2度目の学院生活～"),
    @(27, "さわらないで小手指くん"),
    @(28, "ひゃくえむ。"),
    @(29, "屋根の下のアルテミス"),
    @(30, "降り積もれ孤独な死よ"),
    @(31, "ハードワーカー中田"),
    @(32, "FAIRY TAIL 100 YEARS QUEST"),
    @(33, "アルキメデスの大戦"),
    @(34, "ドラハチ"),
    @(35, "食糧人類-Starving Anonymous-"),
    @(36, "幼馴染とはラブコメにならない"),
    @(37, "ナキナギ"),
    @(38, "グラぱらっ！"),
    @(39, "おやすみ ふみさん"),
    @(40, "Destiny Unchain Online 〜吸血鬼少女となって、やがて『赤の魔王』と呼ばれるようになりました〜"),
    @(41, "辺境の薬師、都でSランク冒険者となる～英雄村の少年がチート薬で無自覚無双〜"),
    @(42, "追放された転生王子、『自動製作《オートクラフト》』スキルで領地を爆速で開拓し最強の村を作ってしまう〜最強クラフトスキルで始める、楽々領地開拓スローライフ〜"),
    @(43, "デッドアカウント"),
    @(44, "時々ボソッとロシア語でデレる隣のアーリャさん"),
    @(45, "アオバノバスケ"),
    @(46, "春くらり"),
    @(47, "我間乱 ―修羅―"),
    @(48, "阿武ノーマル"),
    @(49, "Aランクパーティを離脱した俺は、元教え子たちと迷宮深部を目指す。"),
    @(50, "念願の悪役令嬢（ラスボス）の身体を手に入れたぞ！"),
    @(51, "ジュミドロ"),
    @(52, "ハナバス　苔石花江のバスケ論"),
    @(53, "普通の本はありません！"),
    @(54, "ダメスキル【自動機能】が覚醒しました～あれ、ギルドのスカウトの皆さん、俺を「いらない」って言ってませんでした？～"),
    @(55, "いじめるヤバイ奴"),
    @(56, "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！"),
    @(57, "限界集落を脱村した錬金術士、都会で`"最強`"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～"),
    @(58, "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～"),
    @(59, "なれの果ての僕ら"),
    @(60, "可愛いだけじゃない式守さん"),
    @(61, "君が監督！"),
    @(62, "デスティニーラバーズ"),
    @(63, "東京卍リベンジャーズ～場地圭介からの手紙～"),
    @(64, "劣等人の魔剣使い　スキルボードを駆使して最強に至る"),
    @(65, "インフェクション"),
    @(66, "卒業アルバムの彼女たち"),
    @(67, "不遇職【鑑定士】が実は最強だった～奈落で鍛えた最強の【神眼】で無双する～"),
    @(68, "この世界がいずれ滅ぶことを、俺だけが知っている～モンスターが現れた世界で、死に戻りレベルアップ～"),
    @(69, "MYS"),
    @(70, "恋ニ非ズ"),
    @(71, "母という呪縛 娘という牢獄"),
    @(72, "触手魔術師の成り上がり"),
    @(73, "追放されなかった男　～二度目の人生は土下座から始まりました～"),
    @(74, "リスナーに騙されてダンジョンの最下層から脱出RTAすることになった"),
    @(75, "東京ネオンスキャンダル"),
    @(76, "ぼくのアデリア"),
    @(77, "ストーカー行為がバレて人生終了男"),
    @(78, "シャングリラ・フロンティア～クソゲーハンター、神ゲーに挑まんとす～"),
    @(79, "はっちぽっちぱんち"),
    @(80, "ウイニング パス"),
    @(81, "GALAXIAS"),
    @(82, "正体不明と恐怖"),
    @(83, "人間消失"),
    @(84, "ともだちづくり"),
    @(85, "最弱な僕は＜壁抜けバグ＞で成り上がる～壁をすり抜けたら、初回クリア報酬を無限回収できました！～"),
    @(86, "お願い、脱がシて。"),
    @(87, "勇者と呼ばれた後に　―そして無双男は家族を創る―"),
    @(88, "地獄の業火で焼かれ続けた少年。最強の炎使いとなって復活する。"),
    @(89, "奴隷転生～その奴隷、最強の元王子につき～"),
    @(90, "東大リベンジャーズ"),
    @(91, "中華一番！極"),
    @(92, "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～"),
    @(93, "五輪の女神さま 〜なでしこ寮のメダルごはん〜"),
    @(94, "金田一少年の事件簿外伝 犯人たちの事件簿"),
    @(95, "田んぼで拾った女騎士、田舎で俺の嫁だと思われている"),
    @(96, "彼女、お借りします"),
    @(97, "せいぶつ部の田辺くん"),
    @(98, "不遇職『鍛冶師』だけど最強です ～気づけば何でも作れるようになっていた男ののんびりスローライフ～"),
    @(99, "インフルエンサーにストーカーされています"),
    @(100, "ヒロインは絶望しました。")
)

foreach ($entry in $titles) {
    $rank = $entry[0]
    $title = $entry[1]
    $row = $rank + 1
    $newSheet.Cells.Item($row, 1).Value = $rank
    $newSheet.Cells.Item($row, 2).Value = $title
}
